# Generate Report for Handback
#
# - Status moves from "Ready for handoff" -> "Handed back: in sync with en-US"
#   on the Overview sheet and on each language sheet.
# - Each language sheet gets its "Latest Target File" (F) and
#   "Latest Handback File" (G) columns populated with hyperlinks mirroring
#   the existing "Source File Name" (A) / "Latest Handoff File" (D) links
#   (the handback went out in sync, so target/handback == source/handoff).
# - "Latest Handback DateTime" (H) is stamped with the real handback time
#   instead of the zero-date sentinel. zh-cn and de-de completed at
#   different times so they get different stamps.

$wb = $excel.ActiveWorkbook

function Set-HyperlinkLook($rng) {
    $rng.Style = "HyperLink"
    $rng.Font.Underline = 2
    $rng.Font.Color = 15570276
}

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet: just the status text changes ---
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("B2").Value = $newStatus
$ov.Range("C2").Value = $newStatus
$ov.Range("B3").Value = $newStatus
$ov.Range("C3").Value = $newStatus

# --- zh-cn sheet ---
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("C2").Value = $newStatus
$wsZh.Range("C3").Value = $newStatus

$wsZh.Hyperlinks.Add($wsZh.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/17419e685061977e6a2e84d166e77cff81a12351/e2e/4ee676a3-847d-4da1-ac1e-991f35c7b05f.md", "", "", "4ee676a3-847d-4da1-ac1e-991f35c7b05f.md")
Set-HyperlinkLook $wsZh.Range("F2")

$wsZh.Hyperlinks.Add($wsZh.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a70a85bdc5226d8248e8b0e4f81958fc7b6b4f8f/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/4ee676a3-847d-4da1-ac1e-991f35c7b05f.a8f56980a19e1c46fcc297d63a076c161ed3dc84.zh-cn.xlf", "", "", "4ee676a3-847d-4da1-ac1e-991f35c7b05f.a8f56980a19e1c46fcc297d63a076c161ed3dc84.zh-cn.xlf")
Set-HyperlinkLook $wsZh.Range("G2")

$wsZh.Hyperlinks.Add($wsZh.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/17419e685061977e6a2e84d166e77cff81a12351/e2e/8eabed34-5d2c-4481-a4f6-63c267a3d221.md", "", "", "8eabed34-5d2c-4481-a4f6-63c267a3d221.md")
Set-HyperlinkLook $wsZh.Range("F3")

$wsZh.Hyperlinks.Add($wsZh.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a70a85bdc5226d8248e8b0e4f81958fc7b6b4f8f/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/8eabed34-5d2c-4481-a4f6-63c267a3d221.090b692fd9697d1e9717a86c0bbacea8e20c5e2a.zh-cn.xlf", "", "", "8eabed34-5d2c-4481-a4f6-63c267a3d221.090b692fd9697d1e9717a86c0bbacea8e20c5e2a.zh-cn.xlf")
Set-HyperlinkLook $wsZh.Range("G3")

$wsZh.Range("H2").Value = "2016-03-25 01:25:25"
$wsZh.Range("H3").Value = "2016-03-25 01:25:25"

# --- de-de sheet ---
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value = $newStatus
$wsDe.Range("C3").Value = $newStatus

$wsDe.Hyperlinks.Add($wsDe.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/17419e685061977e6a2e84d166e77cff81a12351/e2e/4ee676a3-847d-4da1-ac1e-991f35c7b05f.md", "", "", "4ee676a3-847d-4da1-ac1e-991f35c7b05f.md")
Set-HyperlinkLook $wsDe.Range("F2")

$wsDe.Hyperlinks.Add($wsDe.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/99333ccf5cc599f3ea9e94806f1804ef26cf9738/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/4ee676a3-847d-4da1-ac1e-991f35c7b05f.a8f56980a19e1c46fcc297d63a076c161ed3dc84.de-de.xlf", "", "", "4ee676a3-847d-4da1-ac1e-991f35c7b05f.a8f56980a19e1c46fcc297d63a076c161ed3dc84.de-de.xlf")
Set-HyperlinkLook $wsDe.Range("G2")

$wsDe.Hyperlinks.Add($wsDe.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/17419e685061977e6a2e84d166e77cff81a12351/e2e/8eabed34-5d2c-4481-a4f6-63c267a3d221.md", "", "", "8eabed34-5d2c-4481-a4f6-63c267a3d221.md")
Set-HyperlinkLook $wsDe.Range("F3")

$wsDe.Hyperlinks.Add($wsDe.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/99333ccf5cc599f3ea9e94806f1804ef26cf9738/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/8eabed34-5d2c-4481-a4f6-63c267a3d221.090b692fd9697d1e9717a86c0bbacea8e20c5e2a.de-de.xlf", "", "", "8eabed34-5d2c-4481-a4f6-63c267a3d221.090b692fd9697d1e9717a86c0bbacea8e20c5e2a.de-de.xlf")
Set-HyperlinkLook $wsDe.Range("G3")

$wsDe.Range("H2").Value = "2016-03-25 01:25:33"
$wsDe.Range("H3").Value = "2016-03-25 01:25:33"

Write-Output "Handback report generated."
